# Update Consumption_Actual.xlsx sheet1 with refreshed data (retraining models)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-42 (column A = consumption value, column B = timestamp serial)
$data = @(
    @(2, 5184, 45918),
    @(3, 5139, 45918.01041666666),
    @(4, 5137, 45918.02083333334),
    @(5, 5095, 45918.03125),
    @(6, 4961, 45918.04166666666),
    @(7, 5014, 45918.05208333334),
    @(8, 4980, 45918.0625),
    @(9, 4962, 45918.07291666666),
    @(10, 4860, 45918.08333333334),
    @(11, 4926, 45918.09375),
    @(12, 4970, 45918.10416666666),
    @(13, 4926, 45918.11458333334),
    @(14, 4970, 45918.125),
    @(15, 4936, 45918.13541666666),
    @(16, 5016, 45918.14583333334),
    @(17, 4971, 45918.15625),
    @(18, 5096, 45918.16666666666),
    @(19, 5138, 45918.17708333334),
    @(20, 5151, 45918.1875),
    @(21, 5247, 45918.19791666666),
    @(22, 5446, 45918.20833333334),
    @(23, 5551, 45918.21875),
    @(24, 5678, 45918.22916666666),
    @(25, 5764, 45918.23958333334),
    @(26, 5974, 45918.25),
    @(27, 6151, 45918.26041666666),
    @(28, 6227, 45918.27083333334),
    @(29, 6263, 45918.28125),
    @(30, 6255, 45918.29166666666),
    @(31, 6172, 45918.30208333334),
    @(32, 6178, 45918.3125),
    @(33, 6171, 45918.32291666666),
    @(34, 5883, 45918.33333333334),
    @(35, 5818, 45918.34375),
    @(36, 5783, 45918.35416666666),
    @(37, 5616, 45918.36458333334),
    @(38, 5324, 45918.375),
    @(39, 5195, 45918.38541666666),
    @(40, 5135, 45918.39583333334),
    @(41, 5055, 45918.40625),
    @(42, 4914, 45918.41666666666)
)

# The existing timestamp format used for column B (rows 2-40)
$dateFormat = $ws.Cells.Item(40, 2).NumberFormat

foreach ($entry in $data) {
    $row = $entry[0]
    $aVal = $entry[1]
    $bVal = $entry[2]

    $ws.Cells.Item($row, 1).Value = $aVal
    $ws.Cells.Item($row, 2).Value = $bVal

    if ($row -gt 40) {
        $ws.Cells.Item($row, 2).NumberFormat = $dateFormat
    }
}
